# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" tab (fund-holding detail, same shape as the
# existing quarter tabs) right before the "2022-Q2" tab, and updates the
# "总计" (summary) tab with a new leading row for 2022-Q3, renumbering the
# A-column sequence (0..5) down through the existing rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New sheet "2022-Q3" - inserted immediately before "2022-Q2"
# ---------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($beforeSheet)
$q3.Name = "2022-Q3"

$headerCols = @("B","C","D","E","F","G","H")
$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headerCols.Length; $i++) {
    $cell = $q3.Range($headerCols[$i] + "1")
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.BorderAround(1)
}

# row 2 - 000646 / 华润元大量化优选混合A
$q3.Range("A2").Value = 0
$a2 = $q3.Range("B2")
$a2.NumberFormat = "@"
$a2.Value = "000646"
$q3.Range("C2").Value = "华润元大量化优选混合A"
$d2 = $q3.Range("D2")
$d2.NumberFormat = "@"
$d2.Value = "1.47"
$e2 = $q3.Range("E2")
$e2.NumberFormat = "@"
$e2.Value = "73.62"
$f2 = $q3.Range("F2")
$f2.NumberFormat = "@"
$f2.Value = "5.61"
$g2 = $q3.Range("G2")
$g2.NumberFormat = "@"
$g2.Value = "0.0825"
$q3.Range("H2").Value = 2

# row 3 - 007827 / 华润元大量化优选混合C
$q3.Range("A3").Value = 1
$b3 = $q3.Range("B3")
$b3.NumberFormat = "@"
$b3.Value = "007827"
$q3.Range("C3").Value = "华润元大量化优选混合C"
$d3 = $q3.Range("D3")
$d3.NumberFormat = "@"
$d3.Value = "0.19"
$e3 = $q3.Range("E3")
$e3.NumberFormat = "@"
$e3.Value = "73.62"
$f3 = $q3.Range("F3")
$f3.NumberFormat = "@"
$f3.Value = "5.61"
$g3 = $q3.Range("G3")
$g3.NumberFormat = "@"
$g3.Value = "0.0107"
$q3.Range("H3").Value = 2

foreach ($col in @("A","B","C","D","E","F","G","H")) {
    $q3.Range($col + "2:" + $col + "3").VerticalAlignment = -4160
}

# ---------------------------------------------------------------------
# 2) "总计" sheet - add the 2022-Q3 row on top and renumber the rest
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$dates = @("2022-Q3","2022-Q2","2022-Q1","2021-Q4","2021-Q3","2021-Q2")
$counts = @(2,2,3,2,1,4)
$values = @(0.09,0.11,0.15,0.36,0.1,0.36)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $i + 2
    $summary.Range("A" + $r).Value = $i
    $summary.Range("B" + $r).Value = $dates[$i]
    $summary.Range("C" + $r).Value = $counts[$i]
    $summary.Range("D" + $r).Value = $values[$i]
}

# Restore the originally-active tab (inserting a sheet makes it active by
# default, but the last tab, "2021-Q2", was the one selected before the edit).
$wb.Worksheets.Item("2021-Q2").Activate()
